# Add a new "Status" column to the asset-management sheet, inserted right
# before the existing "User" column (so the header row becomes:
# ... Nilai Perolehan | Merk -> handled already | Status | User | Dept)
#
# Current header layout (1-indexed columns):
#   A Kode Asset Lama   B Lokasi   C Kategori   D Asset Position  E Merk
#   F Jenis             G Deskripsi  H Serial Number  I Tanggal Perolehan
#   J Umur Ekonomis (Tahun)  K Nilai Perolehan  L User  M Dept
#
# New layout after inserting a "Status" column before column L (User):
#   ... K Nilai Perolehan  L Status  M User  N Dept

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at position 12 (L). This shifts the existing
# "User" (was L) and "Dept" (was M) columns one position to the right,
# and grows the used range from A1:M1 to A1:N1 automatically.
$ws.Columns.Item(12).Insert()

# Populate the header of the newly inserted column.
$ws.Cells.Item(1, 12).Value = "Status"

# Restore/update the view state to match where the user left the cursor
# after making the edit.
$ws.Range("L6").Select()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1

$wb.Save()
